$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the report date references (E2/E3) to the new report date
$ws.Range("E2").Value = "ICER RRMM 2022 report - ICER - 12/19/2022"
$ws.Range("E3").Value = "ICER RRMM 2022 report - ICER - 12/19/2022"

# Merge-select E2:E3 as the active selection, matching the updated selection in the sheet view
$ws.Range("E2:E3").Select()
